$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking price/percent strings
# are preserved exactly (no auto numeric conversion / rounding).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.070.02'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.424.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.98'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.527'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.409.85'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.338'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.851.61'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.828.06'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.421.64'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +7.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.67'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.48'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.02%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.87'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.91'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '592.56'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -8.66%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0937'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.85'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.63'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.03'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.15'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.27'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.70%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0293'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +16.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.38'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.591'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.61'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0504'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.18%  '
